# Update the "已关注/浏览" style counter column (F) on three sheets of the
# "广州-漫展信息" workbook to reflect the freshly generated output.
$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 2794
$ws1.Range("F10").Value = 778
$ws1.Range("F14").Value = 1173
$ws1.Range("F17").Value = 639
$ws1.Range("F21").Value = 709
$ws1.Range("F23").Value = 8240
$ws1.Range("F24").Value = 541
$ws1.Range("F25").Value = 541
$ws1.Range("F30").Value = 219
$ws1.Range("F31").Value = 1697
$ws1.Range("F34").Value = 468
$ws1.Range("F35").Value = 161
$ws1.Range("F38").Value = 172
$ws1.Range("F39").Value = 38
$ws1.Range("F42").Value = 162
$ws1.Range("F45").Value = 30

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F12").Value = 60
$ws2.Range("F15").Value = 48

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 2794
$ws4.Range("F9").Value  = 778
$ws4.Range("F15").Value = 1173
$ws4.Range("F19").Value = 639
$ws4.Range("F25").Value = 8240
$ws4.Range("F27").Value = 541
$ws4.Range("F28").Value = 541
$ws4.Range("F29").Value = 219
$ws4.Range("F30").Value = 1697
$ws4.Range("F32").Value = 468
$ws4.Range("F33").Value = 161
$ws4.Range("F34").Value = 60
$ws4.Range("F35").Value = 60
$ws4.Range("F39").Value = 48
$ws4.Range("F40").Value = 172
$ws4.Range("F41").Value = 38
$ws4.Range("F45").Value = 162

$wb.Save()
